$d = $word.ActiveDocument
$d.Content.Find.Execute("Prueba 1", $true, $false, $false, $false, $false, $true, 1, $false, "Prueba 3", 2)
